$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.387.71'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.826.15'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.70'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9965'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5335'
$ws.Range("E7").Value = '  -1.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3986'
$ws.Range("E8").Value = '  +5.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07545'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.75'
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.104'
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.610'
$ws.Range("E12").Value = '  +3.70%  '
$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9976'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.266'
$ws.Range("E14").Value = '  +1.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.62'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").Value = '1.813.85'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.58'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001068'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06582'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9969'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.029'
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").Value = '28.397.35'
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.16'
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.074'
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.55'
$ws.Range("E26").Value = '  -2.81%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.023.23'
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.390'
$ws.Range("E29").Value = '  +2.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.61'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.108'
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1096'
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.680'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.589'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07387'
$ws.Range("E35").Value = '  +12.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2230'
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.211'
$ws.Range("E37").Value = '  +3.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02304'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.653'
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.31'
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6219'
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.191'
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.405'
$ws.Range("E43").Value = '  -3.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.52'
$ws.Range("E44").Value = '  +0.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.691'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5784'
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '124.95'
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.950'
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.186'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06875'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.29'
$ws.Range("E51").Value = '  -1.55%  '

Write-Host "Cryptos list updated"
